$d = $word.ActiveDocument

# 1) "Web : server, go" -> "Langages : d, postgresql, mysql, python, matlab, c, c++"
$d.Content.Find.Execute(
    "Web : server, go", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Langages : d, postgresql, mysql, python, matlab, c, c++", 2)

# 2) Remove the three paragraphs that followed it:
#    "Autres : marketing, digital marketing, demographic data, barbecues, sql server"
#    "Langages : python, matlab, c, c++"
#    "Soft_Skills : organizational"
$targets = @(
    "Autres : marketing, digital marketing, demographic data, barbecues, sql server",
    "Langages : python, matlab, c, c++",
    "Soft_Skills : organizational"
)
foreach ($target in $targets) {
    $paras = $d.Paragraphs
    for ($i = $paras.Count; $i -ge 1; $i--) {
        $p = $paras.Item($i)
        if ($p.Range.Text.TrimEnd("`r`a") -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}

# 3) "Visualisation : optimization, tableau" -> "Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn"
$d.Content.Find.Execute(
    "Visualisation : optimization, tableau", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", 2)

# 4) "MLOps : hadoop, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit" -> "Visualisation : etl, tableau"
$d.Content.Find.Execute(
    "MLOps : hadoop, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Visualisation : etl, tableau", 2)

# 5) "Maths : algorithms" -> "Machine Learning : hive, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"
$d.Content.Find.Execute(
    "Maths : algorithms", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Machine Learning : hive, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2)

# 6) "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn" -> "Autres : oracle, cassandra"
$d.Content.Find.Execute(
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Autres : oracle, cassandra", 2)
